$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header text changes (force text storage like the rest of the sheet) ---
$ws.Range("J1").Value = "'ΚΩΔΙΚΟΣ ΠΡΟΙΟΝΤΟΣ"
$ws.Range("K1").Value = "'ΟΝΟΜΑ ΠΡΟΙΟΝΤΟΣ"
$ws.Range("L1").Value = "'ΠΟΣΟΤΗΤΑ ΠΡΟΙΟΝΤΟΣ"

# --- Row 2 updates ---
$ws.Range("A2").Value = "'ΜΩΡΑΙΤΗ  ΜΑΡΙΑ ΑΘΑΝΑΣΙΟΣ"
$ws.Range("C2").Value = "'055472412"
$ws.Range("E2").Value = "'2610433413"
$ws.Range("G2").Value = "'26441"
$ws.Range("H2").Value = "'ΠΑΤΡΑ"
$ws.Range("I2").Value = "'ΑΓΙΑΣ ΣΟΦΙΑΣ 15"
$ws.Range("J2").Value = "'GPT-0033"
$ws.Range("K2").Value = "'SAMSUNG TONER CLP325/CLP320/4072 MAGENTA ΣΥΜΒΑΤΟ 1000 ΣΕΛΙΔΕΣ"
$ws.Range("L2").Value = "'15"
$ws.Range("M2").Value = "'1055€"

# --- Row 3 updates ---
$ws.Range("J3").Value = "'GPI-0222"
$ws.Range("K3").Value = "'HP INK No 971XL MAGENTA ΣΥΜΒΑΤΟ 120ml"
$ws.Range("L3").Value = "'98"
$ws.Range("M3").Value = "'"

# --- Row 4 new row ---
$ws.Range("A4").Value = "'"
$ws.Range("B4").Value = "'"
$ws.Range("C4").Value = "'"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "'"
$ws.Range("F4").Value = "'"
$ws.Range("G4").Value = "'"
$ws.Range("H4").Value = "'"
$ws.Range("I4").Value = "'"
$ws.Range("J4").Value = "'GPI-0252"
$ws.Range("K4").Value = "'CANON INK CLI-571XL BLACK ΣΥΜΒΑΤΟ 13ml"
$ws.Range("L4").Value = "'163"
$ws.Range("M4").Value = "'"
